# Pemadanan Data.xlsx - "Coordination Test data Keyword and Listener,
# add Dynamic Date for TC Sanksi"
#
# Update the Kode_PKS test value in the "inquiry" sheet from 01733710 to
# 01733722. The leading apostrophe forces Excel to keep storing the value
# as text (matching the original quote-prefixed numeric-looking string)
# instead of re-interpreting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inquiry")

$ws.Range("A2").Value = "'01733722"
